$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L holds the 2020 data (previous years 2010-2019 are in D:K).
# Row 2 already has an (empty) thin-bottom-border cell under every other
# column in that row, so give L2 the same blank/bordered treatment by
# copying the format from the adjacent K2 cell.
$ws.Range("K2").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null

# Row 3: header cell with year 2020 (bold, bordered, same look as the other
# year headers in D3:K3).
$ws.Range("L3").Value = 2020
$ws.Range("L3").Font.Name = "Times New Roman"
$ws.Range("L3").Font.Size = 9
$ws.Range("L3").Font.Bold = $true
$ws.Range("L3").Borders.Item(9).LineStyle = 1
$ws.Range("L3").Borders.Item(9).Weight = -4138
$ws.Range("L3").Borders.Item(9).ColorIndex = 1

# Data rows 4-9: plain (no border) numeric cells, General format.
$dataRows = @{
    4 = 1004
    5 = 8279
    6 = 1752
    7 = 6527
}
foreach ($r in $dataRows.Keys) {
    $cell = $ws.Range("L$r")
    $cell.Value = $dataRows[$r]
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 9
    $cell.Font.Bold = $false
}

# Rows 8-9: same plain font but number formatted with thousands separator.
$numFmtRows = @{
    8 = 10324
    9 = 4131
}
foreach ($r in $numFmtRows.Keys) {
    $cell = $ws.Range("L$r")
    $cell.Value = $numFmtRows[$r]
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 9
    $cell.Font.Bold = $false
    $cell.NumberFormat = "#,##0"
}

# Row 10: bottom of the table - bordered + thousands separator.
$ws.Range("L10").Value = 6193
$ws.Range("L10").Font.Name = "Times New Roman"
$ws.Range("L10").Font.Size = 9
$ws.Range("L10").Font.Bold = $false
$ws.Range("L10").NumberFormat = "#,##0"
$ws.Range("L10").Borders.Item(9).LineStyle = 1
$ws.Range("L10").Borders.Item(9).Weight = -4138
$ws.Range("L10").Borders.Item(9).ColorIndex = 1

# Selection cursor moved as part of the edit.
$ws.Range("Q11").Select() | Out-Null
